$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.499.70"
$ws.Range("E2").Value = "  +2.05%  "

$ws.Range("D3").Value = "1.856.29"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.88"
$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("E6").Value = "  +0.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9994"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07696"
$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3071"
$ws.Range("E9").Value = "  +0.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.61"
$ws.Range("E10").Value = "  +0.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07779"
$ws.Range("E11").Value = "  -0.28%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.864.33"
$ws.Range("E12").Value = "  +1.66%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.159"
$ws.Range("E13").Value = "  +1.57%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6930"
$ws.Range("E14").Value = "  +1.74%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.08"
$ws.Range("E15").Value = "  +0.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.291"
$ws.Range("E16").Value = "  -2.43%  "

$ws.Range("D17").Value = "29.474.15"
$ws.Range("E17").Value = "  +1.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008337"
$ws.Range("E18").Value = "  +0.53%  "

$ws.Range("D19").Value = "2.100.79"
$ws.Range("E19").Value = "  +1.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.97"
$ws.Range("E20").Value = "  -2.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.73"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.616"
$ws.Range("E23").Value = "  +1.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9997"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1496"
$ws.Range("E25").Value = "  +1.76%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.889"
$ws.Range("E26").Value = "  +0.99%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.69"
$ws.Range("E27").Value = "  -1.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.26"
$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.529"
$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.245"
$ws.Range("E30").Value = "  +0.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.151"
$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.198"
$ws.Range("E32").Value = "  +1.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05099"
$ws.Range("E33").Value = "  -0.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7743"
$ws.Range("E34").Value = "  +1.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.883"
$ws.Range("E35").Value = "  +2.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.149"
$ws.Range("E36").Value = "  +0.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.686"
$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("D38").Value = "1.315.31"
$ws.Range("E38").Value = "  +7.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01874"
$ws.Range("E39").Value = "  +1.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.718"
$ws.Range("E40").Value = "  +0.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9467"
$ws.Range("E41").Value = "  +0.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.16"
$ws.Range("E42").Value = "  -1.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.779"
$ws.Range("E43").Value = "  +1.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.779"
$ws.Range("E45").Value = "  +2.18%  "

$ws.Range("E46").Value = "  +1.90%  "

$ws.Range("D47").Value = "2.001.71"
$ws.Range("E47").Value = "  +1.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5216"
$ws.Range("E48").Value = "  +0.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.789"
$ws.Range("E49").Value = "  +2.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.04"
$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.957"
$ws.Range("E51").Value = "  +0.81%  "
